$d = $word.ActiveDocument

# Change 1: "mine the text we want accordingly." -> "mine the text I want accordingly."
$range1 = $d.Content
$range1.Find.Execute(
    "mine the text we want accordingly.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "mine the text I want accordingly.",
    2)

# Change 2: rewrite the ending of the sentence about presentation time
$range2 = $d.Content
$range2.Find.Execute(
    "and a lack of presentation time. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "and there would be a lack of presentation time to cover this analysis.",
    2)
